$wb = $excel.ActiveWorkbook

# --- Rooms sheet: add a new room row (24) describing the sword tomb ---
$wsRooms = $wb.Worksheets.Item("Rooms")
$wsRooms.Range("A24").Value = "FindSwordRoom1"
$wsRooms.Range("B24").Value = "LootRoom"
$wsRooms.Range("C24").Value = 'You enter a room which has an ominous glow eminating from a center-piece in the middle of the room. On the stone alter a shining knight''s sword lies.  This appears to be a tomb.  The alter reads, "Noble and brave Arturius, son of Theryan, fought and died for the Order of Rosehill."  You need a way to defend yourself, and this sword isn''t being used...'
$wsRooms.Range("D24").Value = "Sword"

# --- Weapons sheet: rebalance the Sword stats (cost/damage/level) ---
$wsWeapons = $wb.Worksheets.Item("Weapons")
$wsWeapons.Range("D5").Value = 30
$wsWeapons.Range("E5").Value = 20
$wsWeapons.Range("G5").Value = 3

# --- Update the active selection on each touched sheet, finishing on the
#     sheet that should remain the active tab (Rooms) ---
$wsArmor = $wb.Worksheets.Item("Armor")
[void]$wsArmor.Range("C26").Select()

[void]$wsWeapons.Range("F5").Select()

[void]$wsRooms.Range("C19").Select()
